$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Remove columns N:O entirely (old "Unnamed: 12"/"Unnamed: 13" + shifted Population type)
$ws.Range("N1:O23").Clear()

# Clear old data rows (2:23) contents; header row A1:M1 keeps its existing bold/border style
$ws.Range("A2:M23").ClearContents()

function ColNum([string]$letter) {
    return ([int][char]$letter) - ([int][char]"A") + 1
}

$table = @(
    @{ Row=1; Cells=@{ "A"="Code Name"; "B"="Display Name"; "C"="Format"; "D"="Timescale"; "E"="Default Value"; "F"="Minimum Value"; "G"="Maximum Value"; "H"="Function"; "I"="Targetable"; "J"="Calibrate"; "K"="Databook Page"; "L"="Databook Order"; "M"="Population type" } },
    @{ Row=2; Cells=@{ "A"="co2e_emissions"; "B"="Total CO2e emissions"; "C"="number"; "H"="energy+bottled_gas_LPG+refridgerants+liquid_fuel+vehicle_fuel_owned+anaesthetic_gases+waste+inhalers"; "I"="n" } },
    @{ Row=3; Cells=@{ "A"="energy_baseline"; "B"="Energy - baseline"; "I"="n"; "K"="emission_sources" } },
    @{ Row=4; Cells=@{ "A"="energy_mult"; "B"="Energy - multiplier"; "E"=0; "F"=0; "G"=1; "I"="y"; "K"="targeted_pars" } },
    @{ Row=5; Cells=@{ "A"="energy"; "B"="Energy"; "H"="energy_baseline*(1-energy_mult)"; "I"="n"; "M"="facilities" } },
    @{ Row=6; Cells=@{ "A"="bottled_gas_LPG_baseline"; "B"="Bottled gas(LPG) - baseline"; "I"="n"; "K"="emission_sources" } },
    @{ Row=7; Cells=@{ "A"="bottled_gas_LPG_mult"; "B"="Bottled gas(LPG) - multiplier"; "E"=0; "F"=0; "G"=1; "I"="y"; "K"="targeted_pars" } },
    @{ Row=8; Cells=@{ "A"="bottled_gas_LPG"; "B"="Bottled gas(LPG)"; "H"="bottled_gas_LPG_baseline*(1-bottled_gas_LPG_mult)"; "I"="n"; "M"="facilities" } },
    @{ Row=9; Cells=@{ "A"="refridgerants_baseline"; "B"="Refridgerants - baseline"; "I"="n"; "K"="emission_sources" } },
    @{ Row=10; Cells=@{ "A"="refridgerants_mult"; "B"="Refridgerants - multiplier"; "E"=0; "F"=0; "G"=1; "I"="y"; "K"="targeted_pars" } },
    @{ Row=11; Cells=@{ "A"="refridgerants"; "B"="Refridgerants"; "H"="refridgerants_baseline*(1-refridgerants_mult)"; "I"="n"; "M"="facilities" } },
    @{ Row=12; Cells=@{ "A"="liquid_fuel_baseline"; "B"="Liquid fuel(Petrol or Diesel) - baseline"; "I"="n"; "K"="emission_sources" } },
    @{ Row=13; Cells=@{ "A"="liquid_fuel_mult"; "B"="Liquid fuel(Petrol or Diesel) - multiplier"; "E"=0; "F"=0; "G"=1; "I"="y"; "K"="targeted_pars" } },
    @{ Row=14; Cells=@{ "A"="liquid_fuel"; "B"="Liquid fuel(Petrol or Diesel)"; "H"="liquid_fuel_baseline*(1-liquid_fuel_mult)"; "I"="n"; "M"="facilities" } },
    @{ Row=15; Cells=@{ "A"="vehicle_fuel_owned_baseline"; "B"="Vehicle Fuel (Owned Vehicles) - baseline"; "I"="n"; "K"="emission_sources" } },
    @{ Row=16; Cells=@{ "A"="vehicle_fuel_owned_mult"; "B"="Vehicle Fuel (Owned Vehicles) - multiplier"; "E"=0; "F"=0; "G"=1; "I"="y"; "K"="targeted_pars" } },
    @{ Row=17; Cells=@{ "A"="vehicle_fuel_owned"; "B"="Vehicle Fuel (Owned Vehicles)"; "H"="vehicle_fuel_owned_baseline*(1-vehicle_fuel_owned_mult)"; "I"="n"; "M"="facilities" } },
    @{ Row=18; Cells=@{ "A"="anaesthetic_gases_baseline"; "B"="Anaesthetic Gases - baseline"; "I"="n"; "K"="emission_sources" } },
    @{ Row=19; Cells=@{ "A"="anaesthetic_gases_mult"; "B"="Anaesthetic Gases - multiplier"; "E"=0; "F"=0; "G"=1; "I"="y"; "K"="targeted_pars" } },
    @{ Row=20; Cells=@{ "A"="anaesthetic_gases"; "B"="Anaesthetic Gases"; "H"="anaesthetic_gases_baseline*(1-anaesthetic_gases_mult)"; "I"="n"; "M"="facilities" } },
    @{ Row=21; Cells=@{ "A"="waste_baseline"; "B"="Waste - baseline"; "I"="n"; "K"="emission_sources" } },
    @{ Row=22; Cells=@{ "A"="waste_mult"; "B"="Waste - multiplier"; "E"=0; "F"=0; "G"=1; "I"="y"; "K"="targeted_pars" } },
    @{ Row=23; Cells=@{ "A"="waste"; "B"="Waste"; "H"="waste_baseline*(1-waste_mult)"; "I"="n"; "M"="facilities" } },
    @{ Row=24; Cells=@{ "A"="inhalers_baseline"; "B"="Inhalers - baseline"; "I"="n"; "K"="emission_sources" } },
    @{ Row=25; Cells=@{ "A"="inhalers_mult"; "B"="Inhalers - multiplier"; "E"=0; "F"=0; "G"=1; "I"="y"; "K"="targeted_pars" } },
    @{ Row=26; Cells=@{ "A"="inhalers"; "B"="Inhalers"; "H"="inhalers_baseline*(1-inhalers_mult)"; "I"="n"; "M"="facilities" } }
)

foreach ($rowdef in $table) {
    $r = $rowdef.Row
    foreach ($col in $rowdef.Cells.Keys) {
        $c = ColNum $col
        $ws.Cells.Item($r, $c).Value = $rowdef.Cells[$col]
    }
}

Write-Output ("New UsedRange: " + $ws.UsedRange.Address())
